$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.315.47'
$ws.Range('E2').Value = '  +4.83%  '
$ws.Range('D3').Value = '1.716.83'
$ws.Range('E3').Value = '  +4.24%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9968'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.35'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9973'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4727'
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2650'
$ws.Range('E8').Value = '  +3.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06232'
$ws.Range('E9').Value = '  +2.12%  '
$ws.Range('D10').Value = '1.700.38'
$ws.Range('E10').Value = '  +3.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07088'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.29'
$ws.Range('E12').Value = '  +6.70%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.5927'
$ws.Range('E13').Value = '  +3.17%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.429'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.31'
$ws.Range('E15').Value = '  +3.88%  '
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9976'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Value = '26.331.75'
$ws.Range('E18').Value = '  +4.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006814'
$ws.Range('E19').Value = '  +2.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.63'
$ws.Range('E20').Value = '  +2.83%  '
$ws.Range('D21').Value = '1.905.83'
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.598'
$ws.Range('E22').Value = '  +6.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.868'
$ws.Range('E23').Value = '  +5.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.359'
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '135.71'
$ws.Range('E25').Value = '  +1.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.19'
$ws.Range('E26').Value = '  +1.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.405'
$ws.Range('E27').Value = '  +1.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.763'
$ws.Range('E28').Value = '  +7.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.69'
$ws.Range('E29').Value = '  +3.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.028'
$ws.Range('E30').Value = '  +3.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.706'
$ws.Range('E31').Value = '  +5.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07773'
$ws.Range('E32').Value = '  +2.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04425'
$ws.Range('E33').Value = '  +4.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.616'
$ws.Range('E34').Value = '  +1.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6248'
$ws.Range('E35').Value = '  +5.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9763'
$ws.Range('E36').Value = '  +4.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9224'
$ws.Range('E37').Value = '  +7.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '110.90'
$ws.Range('E38').Value = '  +12.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.412'
$ws.Range('E39').Value = '  -6.41%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.921'
$ws.Range('E40').Value = '  +7.78%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9992'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01477'
$ws.Range('E42').Value = '  -0.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3820'
$ws.Range('E43').Value = '  +3.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.135'
$ws.Range('E44').Value = '  +11.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1144'
$ws.Range('E45').Value = '  +4.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.273'
$ws.Range('E46').Value = '  +3.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05304'
$ws.Range('E47').Value = '  +1.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.81'
$ws.Range('E48').Value = '  +6.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.691'
$ws.Range('E49').Value = '  +7.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.228'
$ws.Range('E50').Value = '  +2.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3394'
$ws.Range('E51').Value = '  +3.71%  '
